$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New minimal phrase list (header + 5 tongue-twisters)
$ws.Range("A1").Value = "phrase"
$ws.Range("A2").Value = "She sells sea shells by the sea shore"
$ws.Range("A3").Value = "Peter Piper picked a peck of pickled peppers"
$ws.Range("A4").Value = "How much wood would a wood chuck chuck if a wood chuck could chuck wood?"
$ws.Range("A5").Value = "How can a clam cram into a clean cream can?"
$ws.Range("A6").Value = "I saw Susie sitting in a shoeshine shop"

# Remove the now-unused rows that held the old (longer) phrase list
$ws.Range("A7:A11").EntireRow.Delete()

# Move active selection to B2 (matches author's saved selection state)
$ws.Range("B2").Select()
